$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("case study", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ins = $d.Range($r.End, $r.End)
$ins.InsertBefore(" first")
$fmtRange = $d.Range($r.End, $r.End + 6)
$d.Bookmarks.Add("tmpmark", $fmtRange)
$d.Bookmarks("tmpmark").Delete()
